$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (the data currently ends at row 103)
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row

$newRows = @(
    @("2023-12-08 12:07:39", 0.0004),
    @("2023-12-08 12:08:07", 0.0024),
    @("2023-12-08 12:08:42", 0.0026),
    @("2023-12-08 12:08:49", 0.0002)
)

foreach ($entry in $newRows) {
    $lastRow = $lastRow + 1
    $ws.Cells.Item($lastRow, 1).Value = $entry[0]
    $ws.Cells.Item($lastRow, 2).Value = $entry[1]
}
